$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for every data row
# (2..216) from 45184 to 45186.
for ($r = 2; $r -le 216; $r++) {
    $cell = $ws.Range("C$r")
    $cur = $cell.Value2
    if ($cur -eq 45184) {
        $cell.Value = 45186
    }
}

# Add the friendly display text (the "Beteckning" in column A) as the
# second argument of every HYPERLINK() formula in columns S, T, U, V, W,
# X and Y.
$cols = @("S", "T", "U", "V", "W", "X", "Y")
for ($r = 2; $r -le 216; $r++) {
    $a = $ws.Range("A$r").Value2
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula
        if ($f -ne $null -and $f -ne "" -and $f.EndsWith(')') -and -not $f.Contains(', "' + $a + '")')) {
            $cell.Formula = $f.Substring(0, $f.Length - 1) + ', "' + $a + '")'
        }
    }
}
